# Add New Akuntansi Setup Up To TC 9
# Updates the URL cell (D2) on the report-header sheet, normalizes the
# protection-only formatting used across row 2 (dropping the now-unused
# border flag picked up by the old template), and moves the active
# selection, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the URL reported in D2 from the old host to the new host.
$ws.Range("D2").Value = "http://192.168.168.107/"

# Re-assert the "unlocked" protection state on the plain data cells of row 2.
# Doing this per-cell lets the engine re-intern each cell onto the minimal
# matching cell format, dropping the stray applied-border flag that the
# unused/duplicate style carried.
$ws.Range("A2").Locked = $false
$ws.Range("B2").Locked = $false
$ws.Range("C2").Locked = $false
$ws.Range("E2").Locked = $false
$ws.Range("I2").Locked = $false
$ws.Range("K2").Locked = $false
$ws.Range("L2").Locked = $false
$ws.Range("M2").Locked = $false
$ws.Range("U2").Locked = $false
$ws.Range("V2").Locked = $false
$ws.Range("X2").Locked = $false
$ws.Range("Y2").Locked = $false

# Move the active selection/cursor as recorded in the saved view state.
$ws.Range("E10").Select()
